function Set-DateCell {
    param($ws, $addr, $text)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 66 updates ---
$ws.Range("B66").Value = 259.09
$ws.Range("G66").Value = 0

# --- Move existing row 87 ("2023-09-28" ...) data down to row 99 ---
Set-DateCell $ws "A99" "2023-09-28"
$ws.Range("B99").Value = 0
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 100
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0

# --- Move existing row 88 ("2023-11-24" ...) data down to row 100 ---
Set-DateCell $ws "A100" "2023-11-24"
$ws.Range("B100").Value = 0
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 26.01
$ws.Range("E100").Value = 0
$ws.Range("F100").Value = 51.52
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0

# --- Update row 87 with new data ---
Set-DateCell $ws "A87" "2023-11-28"
$ws.Range("B87").Value = 270.57
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 102.53
$ws.Range("E87").Value = 170
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0

# --- New rows 88-98 ---
# Row 88: 2023-11-29
Set-DateCell $ws "A88" "2023-11-29"
$ws.Range("B88").Value = 103.48
$ws.Range("C88").Value = 50
$ws.Range("D88").Value = 224.45
$ws.Range("E88").Value = 45
$ws.Range("F88").Value = 10.7
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0

# Row 89: 2023-11-30
Set-DateCell $ws "A89" "2023-11-30"
$ws.Range("B89").Value = 432.11
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("E89").Value = 0
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0

# Row 90: 2023-12-01
Set-DateCell $ws "A90" "2023-12-01"
$ws.Range("B90").Value = 144.35
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 77.77000000000001
$ws.Range("E90").Value = 600
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0

# Row 91: 2023-12-02
Set-DateCell $ws "A91" "2023-12-02"
$ws.Range("B91").Value = 437.15
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 46.83
$ws.Range("E91").Value = 1000
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0

# Row 92: 2023-12-03
Set-DateCell $ws "A92" "2023-12-03"
$ws.Range("B92").Value = 31.4
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 41.31
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0

# Row 93: 2023-12-04
Set-DateCell $ws "A93" "2023-12-04"
$ws.Range("B93").Value = 73.90000000000001
$ws.Range("C93").Value = 103.48
$ws.Range("D93").Value = 200
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 103.48
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0

# Row 94: 2023-12-05
Set-DateCell $ws "A94" "2023-12-05"
$ws.Range("B94").Value = 122.53
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 50
$ws.Range("E94").Value = 141
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0

# Row 95: 2023-12-06
Set-DateCell $ws "A95" "2023-12-06"
$ws.Range("B95").Value = 76.52000000000001
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 404.56
$ws.Range("E95").Value = 0
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0

# Row 96: 2023-12-07
Set-DateCell $ws "A96" "2023-12-07"
$ws.Range("B96").Value = 630.76
$ws.Range("C96").Value = 113.23
$ws.Range("D96").Value = 35
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0

# Row 97: 2023-12-08
Set-DateCell $ws "A97" "2023-12-08"
$ws.Range("B97").Value = 71.52000000000001
$ws.Range("C97").Value = 307.09
$ws.Range("D97").Value = 300
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0

# Row 98: 2023-12-09
Set-DateCell $ws "A98" "2023-12-09"
$ws.Range("B98").Value = 1204.56
$ws.Range("C98").Value = 77.02
$ws.Range("D98").Value = 51.52
$ws.Range("E98").Value = 510
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0

Write-Host "Done applying edits."
